$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.012.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.561.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0596"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.783.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.563.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("E14").Value = "  -0.07%  "

$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.030.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("E20").Value = "  +1.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("E24").Value = "  -0.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("E28").Value = "  +1.63%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("E31").Value = "  +3.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.421.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.13%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.57%  "

$ws.Range("E37").Value = "  +2.75%  "

$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("E39").Value = "  +2.00%  "

$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("E42").Value = "  +0.33%  "

$ws.Range("E43").Value = "  -0.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("E46").Value = "  -1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.697.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("E49").Value = "  +2.99%  "

$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
